# Update the "Scenario" sheet's Options/Options values table to reflect the
# new pydantic-settings-based MDA options representation: the individual
# warm_start/tolerance/over_relaxation_factor/max_mda_iter rows collapse into
# a single "main_mda_settings" row holding a JSON blob, and "name" moves up
# to take the vacated row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# Row 3: was warm_start / False -> becomes name / my_test_scenario
$ws.Range("L3").Value = "name"
$ws.Range("M3").Value = "my_test_scenario"

# Row 4: was name / my_test_scenario -> becomes main_mda_settings / {json}
$ws.Range("L4").Value = "main_mda_settings"
$ws.Range("M4").Value = '{"max_mda_iter": 20, "warm_start": True, "tolerance": 1e-5, "over_relaxation_factor": 1.2}'

# Rows 5-7 (tolerance, over_relaxation_factor, max_mda_iter) are no longer
# needed - their info now lives in the main_mda_settings JSON blob.
$ws.Range("L5:M7").ClearContents()

$ws.Range("M4").Select() | Out-Null
